# Generate Report for handoff
# The 56bec4ee-7eef-4f4f-8950-234bdbffa32a.md file has been newly handed off again;
# update its status (Overview + both locale sheets) to "Ready for handoff" and stamp
# the new "Latest Handoff Datetime" per locale.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the Status columns for the 56bec4ee... row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: update Status + Latest Handoff Datetime for row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-22 02:30:52"

# --- de-de sheet: update Status + Latest Handoff Datetime for row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-22 02:31:06"
